$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 525.5417
$ws.Range("I53").Value = 615.3
$ws.Range("J53").Value = 461.42856
$ws.Range("K53").Value = 615.3
$ws.Range("L53").Value = 461.42856
$ws.Range("M53").Value = 21.70000000000005
$ws.Range("N53").Value = -1735.42856
# Row 86
$ws.Range("H86").Value = 4193.722
$ws.Range("I86").Value = 2082.6667
$ws.Range("J86").Value = 8415.833000000001
$ws.Range("K86").Value = 2082.6667
$ws.Range("L86").Value = 8415.833000000001
$ws.Range("M86").Value = -959.6667000000002
$ws.Range("N86").Value = -10661.833
# Row 89
$ws.Range("H89").Value = 4193.722
$ws.Range("I89").Value = 2082.6667
$ws.Range("J89").Value = 8415.833000000001
$ws.Range("K89").Value = 10413.3335
$ws.Range("L89").Value = 42079.165
$ws.Range("M89").Value = -4797.333500000001
$ws.Range("N89").Value = -53311.165
# Row 106
$ws.Range("H106").Value = 4799.857
$ws.Range("I106").Value = 5119.8
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 5119.8
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -4488.8
$ws.Range("N106").Value = -5262
# Row 132
$ws.Range("H132").Value = 7359150
$ws.Range("I132").Value = 9266863
$ws.Range("J132").Value = 828.5714
$ws.Range("K132").Value = 27800589
$ws.Range("L132").Value = 2485.7142
$ws.Range("M132").Value = -27798059
$ws.Range("N132").Value = -7545.7142
# Row 137
$ws.Range("H137").Value = 973.2646999999999
$ws.Range("I137").Value = 949.2373
$ws.Range("J137").Value = 1130.7778
$ws.Range("K137").Value = 2847.7119
$ws.Range("L137").Value = 3392.3334
$ws.Range("M137").Value = -297.7119000000002
$ws.Range("N137").Value = -8492.3334
# Row 138
$ws.Range("H138").Value = 1566.9767
$ws.Range("I138").Value = 1162.6571
$ws.Range("J138").Value = 3335.875
$ws.Range("K138").Value = 3487.9713
$ws.Range("L138").Value = 10007.625
$ws.Range("M138").Value = 1652.0287
$ws.Range("N138").Value = -20287.625

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 18
$ws.Range("H18").Value = 49800
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 49800
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 49800
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -50444
# Row 32
$ws.Range("H32").Value = 19945.783
$ws.Range("I32").Value = 4444.9185
$ws.Range("J32").Value = 141134.36
$ws.Range("K32").Value = 4444.9185
$ws.Range("L32").Value = 141134.36
$ws.Range("M32").Value = -4157.9185
$ws.Range("N32").Value = -141708.36
# Row 61
$ws.Range("H61").Value = 1378.6571
$ws.Range("I61").Value = 1126.6451
$ws.Range("J61").Value = 3331.75
$ws.Range("K61").Value = 1126.6451
$ws.Range("L61").Value = 3331.75
$ws.Range("M61").Value = -914.6451
$ws.Range("N61").Value = -3755.75
# Row 74
$ws.Range("H74").Value = 422.65625
$ws.Range("I74").Value = 422.65625
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 422.65625
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 451.34375
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 422.65625
$ws.Range("I77").Value = 422.65625
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2113.28125
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 2254.71875
$ws.Range("N77").ClearContents()
# Row 136
$ws.Range("H136").Value = 1378.6571
$ws.Range("I136").Value = 1126.6451
$ws.Range("J136").Value = 3331.75
$ws.Range("K136").Value = 3379.9353
$ws.Range("L136").Value = 9995.25
$ws.Range("M136").Value = -829.9353000000001
$ws.Range("N136").Value = -15095.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 457.36
$ws.Range("I94").Value = 392.09525
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 392.09525
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 58.90474999999998
$ws.Range("N94").Value = -1702
# Row 134
$ws.Range("H134").Value = 2319.2744
$ws.Range("I134").Value = 2030.8959
$ws.Range("J134").Value = 6933.3335
$ws.Range("K134").Value = 6092.6877
$ws.Range("L134").Value = 20800.0005
$ws.Range("M134").Value = -3557.6877
$ws.Range("N134").Value = -25870.0005

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 29
$ws.Range("H29").Value = 29000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 29000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 29000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -29586
# Row 31
$ws.Range("H31").Value = 30807.3
$ws.Range("I31").Value = 1765.7693
$ws.Range("J31").Value = 62268.957
$ws.Range("K31").Value = 1765.7693
$ws.Range("L31").Value = 62268.957
$ws.Range("M31").Value = -1470.7693
$ws.Range("N31").Value = -62858.957
# Row 34
$ws.Range("H34").Value = 30807.3
$ws.Range("I34").Value = 1765.7693
$ws.Range("J34").Value = 62268.957
$ws.Range("K34").Value = 1765.7693
$ws.Range("L34").Value = 62268.957
$ws.Range("M34").Value = -1563.7693
$ws.Range("N34").Value = -62672.957
# Row 58
$ws.Range("H58").Value = 1303.4546
$ws.Range("I58").Value = 1137.1052
$ws.Range("J58").Value = 2357
$ws.Range("K58").Value = 1137.1052
$ws.Range("L58").Value = 2357
$ws.Range("M58").Value = -934.1052
$ws.Range("N58").Value = -2763
# Row 132
$ws.Range("H132").Value = 3808.1082
$ws.Range("I132").Value = 3695.84
$ws.Range("J132").Value = 4042
$ws.Range("K132").Value = 11087.52
$ws.Range("L132").Value = 12126
$ws.Range("M132").Value = -8557.52
$ws.Range("N132").Value = -17186
# Row 134
$ws.Range("H134").Value = 1019.1842
$ws.Range("I134").Value = 822.96295
$ws.Range("J134").Value = 1500.8182
$ws.Range("K134").Value = 2468.88885
$ws.Range("L134").Value = 4502.4546
$ws.Range("M134").Value = 66.11115000000018
$ws.Range("N134").Value = -9572.454600000001
# Row 136
$ws.Range("H136").Value = 1303.4546
$ws.Range("I136").Value = 1137.1052
$ws.Range("J136").Value = 2357
$ws.Range("K136").Value = 3411.3156
$ws.Range("L136").Value = 7071
$ws.Range("M136").Value = -861.3155999999999
$ws.Range("N136").Value = -12171

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 9077.675999999999
$ws.Range("I131").Value = 5000
$ws.Range("J131").Value = 9135.929
$ws.Range("K131").Value = 15000
$ws.Range("L131").Value = 27407.787
$ws.Range("M131").Value = -9960
$ws.Range("N131").Value = -37487.787

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 802.8182
$ws.Range("I122").Value = 616.5
$ws.Range("J122").Value = 1299.6666
$ws.Range("K122").Value = 1849.5
$ws.Range("L122").Value = 3898.9998
$ws.Range("M122").Value = 600.5
$ws.Range("N122").Value = -8798.9998
# Row 132
$ws.Range("H132").Value = 4095.037
$ws.Range("I132").Value = 4031.1428
$ws.Range("J132").Value = 4318.6665
$ws.Range("K132").Value = 12093.4284
$ws.Range("L132").Value = 12955.9995
$ws.Range("M132").Value = -9563.428400000001
$ws.Range("N132").Value = -18015.9995

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 64654.812
$ws.Range("I40").Value = 201395.8
$ws.Range("J40").Value = 2499.818
$ws.Range("K40").Value = 201395.8
$ws.Range("L40").Value = 2499.818
$ws.Range("M40").Value = -201259.8
$ws.Range("N40").Value = -2771.818
# Row 122
$ws.Range("H122").Value = 2167.6785
$ws.Range("I122").Value = 2043.7391
$ws.Range("J122").Value = 2737.8
$ws.Range("K122").Value = 6131.2173
$ws.Range("L122").Value = 8213.400000000001
$ws.Range("M122").Value = -3681.2173
$ws.Range("N122").Value = -13113.4
# Row 132
$ws.Range("H132").Value = 2549.7083
$ws.Range("I132").Value = 2372.4546
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 7117.3638
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -4587.3638
$ws.Range("N132").Value = -18558.5
# Row 136
$ws.Range("H136").Value = 1071.4717
$ws.Range("I136").Value = 882.913
$ws.Range("J136").Value = 2310.5715
$ws.Range("K136").Value = 2648.739
$ws.Range("L136").Value = 6931.7145
$ws.Range("M136").Value = -98.73900000000003
$ws.Range("N136").Value = -12031.7145

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1889.1311
$ws.Range("I132").Value = 1762.9056
$ws.Range("J132").Value = 2725.375
$ws.Range("K132").Value = 5288.7168
$ws.Range("L132").Value = 8176.125
$ws.Range("M132").Value = -2758.7168
$ws.Range("N132").Value = -13236.125
# Row 136
$ws.Range("H136").Value = 601.6727
$ws.Range("I136").Value = 392.95557
$ws.Range("J136").Value = 1540.9
$ws.Range("K136").Value = 1178.86671
$ws.Range("L136").Value = 4622.700000000001
$ws.Range("M136").Value = 1371.13329
$ws.Range("N136").Value = -9722.700000000001

Write-Output "Applied Aegis_Profits updates"